$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (nombre_aides, montant_total) as published with the
# 2020-08-13 data refresh. Source values are text (as in the original
# sheet), so each cell is forced to Text format before the value is
# written, matching the original file's string-typed numeric columns.
$updates = @(
    @(2,  "162",  "358316.00"),
    @(3,  "890",  "2559892.01"),
    @(4,  "370",  "1399326.69"),
    @(5,  "99",   "431482.09"),
    @(6,  "23",   "144500.00"),
    @(8,  "30",   "60000.00"),
    @(15, "99",   "264152.38"),
    @(16, "414",  "1272472.14"),
    @(18, "46",   "231045.00"),
    @(23, "116",  "468348.00"),
    @(24, "38",   "190233.00"),
    @(26, "18",   "45500.00"),
    @(33, "101",  "270748.00"),
    @(34, "527",  "1639546.26"),
    @(46, "74",   "322164.35"),
    @(47, "41",   "241937.00"),
    @(48, "25",   "173697.00"),
    @(49, "10",   "23850.00"),
    @(56, "679",  "1701943.84"),
    @(57, "3326", "10062234.69"),
    @(58, "1714", "6894260.28"),
    @(59, "582",  "2757640.96"),
    @(60, "119",  "811623.00"),
    @(62, "277",  "658802.79"),
    @(79, "223",  "565326.09"),
    @(80, "856",  "2654635.61"),
    @(81, "322",  "1242468.79"),
    @(82, "109",  "515984.52")
)

foreach ($u in $updates) {
    $row = $u[0]
    $cCell = $ws.Cells.Item($row, 3)
    $cCell.NumberFormat = "@"
    $cCell.Value = $u[1]

    $dCell = $ws.Cells.Item($row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $u[2]
}
